# Applies the "plots for teh gates action and linear models of water quality"
# commit to the FMWT/SMSCG results workbook:
#   - adds a "Dry years" section heading above the existing dry-years block
#     (row 25, in front of the block that starts at what becomes row 26)
#   - appends a brand new "Dry years, scaled covariates" results block
#     (rows 51-68) with count-model and zero-inflation-model coefficient
#     tables, mirroring the layout already used for the other blocks
#   - updates the view's selection to Q55, matching the author's last
#     on-screen selection when they saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-Style($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy()
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# New block at rows 51-68: "Dry years, scaled covariates"
# (Value set first so the shared-string table gets "Dry years, scaled
# covariates" before the shorter "Dry years" heading below - matching
# the order the strings were originally authored in.)
# ---------------------------------------------------------------------

# Row 51 - block title, no explicit style
$ws.Range("A51").Value = "Dry years, scaled covariates"

# ---------------------------------------------------------------------
# New row 25: a short "Dry years" heading inserted above the existing
# "Count model coefficients..." block that used to start at row 26.
# (No explicit style - matches default/no "s" attribute cell.)
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Dry years"

# Row 52 - "Count model coefficients (negbin with log link):"
Copy-Style "A1" "A52"
$ws.Range("A52").Value = "Count model coefficients (negbin with log link):"

# Row 53 - column headers
Copy-Style "A1" "A53"
$ws.Range("B53").Value = "Estimate"
$ws.Range("C53").Value = "Std. Error"
$ws.Range("D53").Value = "z value"
$ws.Range("E53").Value = "Pr(>|z|)"

# Row 54 - (Intercept)
Copy-Style "A1" "A54"
$ws.Range("A54").Value = "(Intercept)"
$ws.Range("B54").Value = 0.4605
$ws.Range("C54").Value = 0.4712
$ws.Range("D54").Value = 0.977
$ws.Range("E54").Value = 0.3285

# Row 55 - Gate Ops
Copy-Style "A1" "A55"
$ws.Range("A55").Value = "Gate Ops"
$ws.Range("B55").Value = -1.3133
$ws.Range("C55").Value = 0.5812
$ws.Range("D55").Value = -2.26
$ws.Range("E55").Value = 0.0238
$ws.Range("F55").Value = "*"

# Row 56 - Day of year
Copy-Style "A1" "A56"
$ws.Range("A56").Value = "Day of year"
$ws.Range("B56").Value = 0.5272
$ws.Range("C56").Value = 0.2363
$ws.Range("D56").Value = 2.231
$ws.Range("E56").Value = 0.0257
$ws.Range("F56").Value = "*"

# Row 57 - Conductivity
Copy-Style "A1" "A57"
$ws.Range("A57").Value = "Conductivity"
$ws.Range("B57").Value = -0.3449
$ws.Range("C57").Value = 0.2227
$ws.Range("D57").Value = -1.549
$ws.Range("E57").Value = 0.1215

# Row 58 - FMWT Index
Copy-Style "A1" "A58"
$ws.Range("A58").Value = "FMWT Index"
$ws.Range("B58").Value = 0.2865
$ws.Range("C58").Value = 0.1562
$ws.Range("D58").Value = 1.834
$ws.Range("E58").Value = 0.0666
$ws.Range("F58").Value = "."

# Row 59 - Log(theta)
Copy-Style "A1" "A59"
$ws.Range("A59").Value = "Log(theta)"
$ws.Range("B59").Value = -0.1285
$ws.Range("C59").Value = 0.4563
$ws.Range("D59").Value = -0.282
$ws.Range("E59").Value = 0.7783

# Row 60 - blank spacer row
Copy-Style "A11" "A60"

# Row 61 - "Zero-inflation model coefficients (binomial with logit link):"
Copy-Style "A1" "A61"
$ws.Range("A61").Value = "Zero-inflation model coefficients (binomial with logit link):"

# Row 62 - column headers
Copy-Style "A1" "A62"
$ws.Range("B62").Value = "Estimate"
$ws.Range("C62").Value = "Std. Error"
$ws.Range("D62").Value = "z value"
$ws.Range("E62").Value = "Pr(>|z|)"

# Row 63 - (Intercept)
Copy-Style "A1" "A63"
$ws.Range("A63").Value = "(Intercept)"
$ws.Range("B63").Value = -65.24
$ws.Range("C63").Value = 83.4
$ws.Range("D63").Value = -0.782
$ws.Range("E63").Value = 0.434

# Row 64 - Operating
Copy-Style "A1" "A64"
$ws.Range("A64").Value = "Operating"
$ws.Range("B64").Value = -189.06
$ws.Range("C64").Value = 233.06
$ws.Range("D64").Value = -0.811
$ws.Range("E64").Value = 0.417

# Row 65 - Day of year
Copy-Style "A1" "A65"
$ws.Range("A65").Value = "Day of year"
$ws.Range("B65").Value = 44.54
$ws.Range("C65").Value = 58.02
$ws.Range("D65").Value = 0.768
$ws.Range("E65").Value = 0.443

# Row 66 - Conductivity
Copy-Style "A1" "A66"
$ws.Range("A66").Value = "Conductivity"
$ws.Range("B66").Value = -52.58
$ws.Range("C66").Value = 64.34
$ws.Range("D66").Value = -0.817
$ws.Range("E66").Value = 0.414

# Row 67 - FMWT Index
Copy-Style "A1" "A67"
$ws.Range("A67").Value = "FMWT Index"
$ws.Range("B67").Value = -225.33
$ws.Range("C67").Value = 277.66
$ws.Range("D67").Value = -0.812
$ws.Range("E67").Value = 0.417

# Row 68 - closing "---" separator
Copy-Style "A21" "A68"
$ws.Range("A68").Value = "---"

# ---------------------------------------------------------------------
# View state: scroll so row 34 is at the top and select Q55, matching
# the author's on-screen state at save time.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q55").Select()
